$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1.28
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0.9
$ws.Range("F2").Value = 2.9
$ws.Range("G2").Value = 3

# Row 3
$ws.Range("C3").Value = 1.02
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = 0.9
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 3
$ws.Range("I3").Value = 0.07000000000000001

# Row 4
$ws.Range("B4").Value = 2.38
$ws.Range("C4").Value = 3
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0.9
$ws.Range("F4").Value = 2.99
$ws.Range("G4").Value = 2.99
